# The commit simplifies the document-wide defaults in word/styles.xml:
# <w:docDefaults> currently spells out every formatting knob explicitly
# (bold=0, italic=0, no shading, widowControl=1, zero indents, jc=left, …)
# even though those are all the OOXML spec's own implicit defaults. The
# change drops the redundant, default-valued properties and keeps only the
# ones that actually carry meaning: the Arial/22half-pt/"en" run defaults
# and the 276-auto line spacing paragraph default.
#
# This is a package-level (docDefaults) edit with no Word-UI/object-model
# surface of its own (Word never exposes docDefaults through Styles/Font/
# ParagraphFormat — those always read/write the *Normal* style's own
# rPr/pPr, not the template-level docDefaults block). The supported COM
# mechanism for this kind of direct, surgical WordprocessingML edit is the
# document's WordOpenXML round-trip: pull the flat-OPC XML, rewrite the
# <w:docDefaults> element, and push it back.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

$pattern = '<w:docDefaults>.*?</w:docDefaults>'

$replacement = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:pPrDefault></w:docDefaults>'

$regex = New-Object System.Text.RegularExpressions.Regex($pattern, [System.Text.RegularExpressions.RegexOptions]::Singleline)

$newXml = $regex.Replace($xml, $replacement, 1)

$d.WordOpenXML = $newXml
